# edit.ps1 -- applies the coinranking "symbol list" refresh captured by the
# GitHub Actions commit "Updated symbol list on Fri Dec 16 16:34:31 UTC 2022".
#
# The sheet stores every data cell as text (inlineStr), including the Price
# column, which holds numeric-looking values. Assigning a plain numeric string
# to Range.Value lets Excel auto-convert it to a real number, so numeric-looking
# values are written with a leading apostrophe (Excel's "force text" quote
# prefix) to keep them stored as text, matching the workbook's existing format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BNB) -- price refresh
$ws.Range("D2").Value = "'243.26"
$ws.Range("D3").Value = "'23.86"
$ws.Range("D4").Value = "'5.758"
$ws.Range("D5").Value = "'0.05845"
$ws.Range("D6").Value = "'3.419"
$ws.Range("D7").Value = "'6.508"
$ws.Range("D8").Value = "'1.325"
$ws.Range("D9").Value = "'0.7991"

# Rows 10-18 -- coin ranking reshuffled down one slot, "One" rises to rank 9
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01252"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1473"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07693"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03301"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03020"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09212"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.576"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001667"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04762"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Rows 19-27 -- price refresh
$ws.Range("D19").Value = "'0.006263"
$ws.Range("D20").Value = "'0.005485"
$ws.Range("D21").Value = "'0.001071"
$ws.Range("D22").Value = "'0.0001503"
$ws.Range("D23").Value = "'3.716"
$ws.Range("D25").Value = "'0.3327"
$ws.Range("D26").Value = "'0.1255"
$ws.Range("D27").Value = "'0.0006282"

# Rows 40-41 -- price refresh
$ws.Range("D40").Value = "'0.04325"
$ws.Range("D41").Value = "'0.007045"

# Rows 42-43 -- BKEXToken / CEJI swap ranks
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1055"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003416"
$ws.Range("E43").Value = "42CEJICEJI"

# Rows 44-46 -- price refresh + best/worst-in-24h label churn
$ws.Range("D44").Value = "'0.008729"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"
$ws.Range("D45").Value = "'0.002469"
$ws.Range("E45").Value = "44ACDXExchangeACXT"
$ws.Range("D46").Value = "'0.00005765"

# Rows 48-51 -- price refresh + best/worst-in-24h label churn
$ws.Range("D48").Value = "'0.9919"
$ws.Range("D49").Value = "'0.1075"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D51").Value = "'0.01012"
